$wb = $excel.ActiveWorkbook

$xlPasteFormats = -4122
$xlEdgeTop = 8
$xlEdgeBottom = 9
$xlEdgeRight = 10
$xlContinuous = 1

$ws1 = $wb.Worksheets.Item("quality_comparison")
$ws2 = $wb.Worksheets.Item("computational_comparison")

# --- Build the two new border styles once (on sheet1 C1/D1), then replicate
# them to every other header cell that needs the same border via
# copy/paste-special-formats, so every cell lands on its target style in a
# single step (avoids generating extra, unused intermediate cell styles).

# "Middle" header cell style: thin border on top & bottom only.
$midRef = $ws1.Range("C1")
$midRef.ClearFormats()
$midRef.Borders.Item($xlEdgeTop).LineStyle = $xlContinuous
$midRef.Borders.Item($xlEdgeBottom).LineStyle = $xlContinuous

# "Right-most" header cell style: thin border on top, right & bottom.
$midRef.Copy()
$rightRef = $ws1.Range("D1")
$rightRef.PasteSpecial($xlPasteFormats)
$rightRef.Borders.Item($xlEdgeRight).LineStyle = $xlContinuous

# Apply the same two styles to the matching header cells on sheet2.
$midRef.Copy()
$ws2.Range("C1").PasteSpecial($xlPasteFormats)

$rightRef.Copy()
$ws2.Range("D1").PasteSpecial($xlPasteFormats)

$midRef.Copy()
$ws2.Range("F1").PasteSpecial($xlPasteFormats)

$rightRef.Copy()
$ws2.Range("G1").PasteSpecial($xlPasteFormats)

$excel.CutCopyMode = 0

# --- Rename the "fedcore" column header to "approach" (anonymized) ---
$ws1.Range("C2").Value = "approach"
$ws2.Range("C2").Value = "approach"
$ws2.Range("F2").Value = "approach"

# --- Drop the stray empty inline-string cell at G5 on sheet2 ---
$ws2.Range("G5").ClearContents()
